# Update open/close/high/low/shares_outstanding/fixed_ticker columns (D:I)
# for each data row (2-23) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  D=84;                 E=86.41999816894531; F=91.98999786376952; G=78.95999908447266; H=150192801; I="PCOR"},
    @{Row=3;  D=84;                 E=86.41999816894531; F=91.98999786376952; G=78.95999908447266; H=150192801; I="PCOR"},
    @{Row=4;  D=84;                 E=86.41999816894531; F=91.98999786376952; G=78.95999908447266; H=150192801; I="PCOR"},
    @{Row=5;  D=84;                 E=86.41999816894531; F=91.98999786376952; G=78.95999908447266; H=150192801; I="PCOR"},
    @{Row=6;  D=84;                 E=86.41999816894531; F=91.98999786376952; G=78.95999908447266; H=150192801; I="PCOR"},
    @{Row=7;  D=94.76000213623048; E=103.2799987792969; F=104.4199981689453; G=90.05000305175781; H=150192801; I="PCOR"},
    @{Row=8;  D=89.44999694824219; E=91.4499969482422;  F=105.2900009155273; G=85.20500183105469; H=150192801; I="PCOR"},
    @{Row=9;  D=80.23000335693359; E=62.56000137329102; F=80.81099700927734; G=56;                 H=150192801; I="PCOR"},
    @{Row=10; D=58.36999893188477; E=55.47000122070312; F=62.77000045776367; G=49.90999984741211; H=150192801; I="PCOR"},
    @{Row=11; D=45.66999816894531; E=51.70000076293945; F=52.5;              G=43.45000076293945; H=150192801; I="PCOR"},
    @{Row=12; D=50.02000045776367; E=54.65999984741211; F=57.02999877929688; G=45.45000076293945; H=150192801; I="PCOR"},
    @{Row=13; D=48.18999862670898; E=55.95000076293945; F=58.09000015258789; G=41.79999923706055; H=150192801; I="PCOR"},
    @{Row=14; D=62.11999893188477; E=53.40999984741211; F=62.3380012512207;  G=51.36999893188477; H=150192801; I="PCOR"},
    @{Row=15; D=64.5;               E=75.84999847412109; F=76.25;             G=61.38000106811523; H=150192801; I="PCOR"},
    @{Row=16; D=65.06999969482422; E=61.09000015258789; F=70.27999877929688; G=60.13000106811523; H=150192801; I="PCOR"},
    @{Row=17; D=68.75;              E=71.38999938964844; F=73.44000244140625; G=64.04000091552734; H=150192801; I="PCOR"},
    @{Row=18; D=82.15000152587891; E=68.41999816894531; F=82.63999938964844; G=68.36000061035156; H=150192801; I="PCOR"},
    @{Row=19; D=66.12999725341797; E=71.02999877929688; F=72.48000335693359; G=63.33000183105469; H=150192801; I="PCOR"},
    @{Row=20; D=61.68999862670898; E=65.65000152587891; F=69.27999877929688; G=58.09999847412109; H=150192801; I="PCOR"},
    @{Row=21; D=76.12999725341797; E=79.55999755859375; F=83.06999969482422; G=72.88099670410156; H=150192801; I="PCOR"},
    @{Row=22; D=66.19000244140625; E=64.08999633789062; F=69.98999786376953; G=53.70999908447266; H=150192801; I="PCOR"},
    @{Row=23; D=68.20999908447266; E=71.62999725341797; F=77.88999938964844; G=67.48000335693359; H=150192801; I="PCOR"}
)

foreach ($rowData in $data) {
    $r = $rowData.Row
    $ws.Range("D$r").Value = $rowData.D
    $ws.Range("E$r").Value = $rowData.E
    $ws.Range("F$r").Value = $rowData.F
    $ws.Range("G$r").Value = $rowData.G
    $ws.Range("H$r").Value = $rowData.H
    $ws.Range("I$r").Value = $rowData.I
}
